# Apply updated crypto market data (prices & 1h volume %) scraped on
# Sat May 25 14:32:42 UTC 2024, plus reordering of the Bittensor/
# dogwifhat rows (ranks 40/41) to reflect the new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain-looking numbers as text (e.g. "1.00",
# "601.17"); force text format first so Excel does not coerce them to
# numeric values and drop the formatting (trailing zeros, etc).
$textCells = @("D5","D6","D11","D12","D13","D18","D20","D21","D22","D23","D24","D25","D27","D28","D31","D33","D37","D38","D39","D41","D42","D43","D44","D45","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '69.023.08'
$ws.Range("E2").Value = '  +1.74%  '

# Row 3
$ws.Range("D3").Value = '3.738.73'
$ws.Range("E3").Value = '  +0.60%  '

# Row 4
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").Value = '601.17'
$ws.Range("E5").Value = '  +0.39%  '

# Row 6
$ws.Range("D6").Value = '167.57'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").Value = '3.737.98'
$ws.Range("E7").Value = '  +0.64%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("E10").Value = '  +0.10%  '

# Row 11
$ws.Range("D11").Value = '6.44'
$ws.Range("E11").Value = '  +3.60%  '

# Row 12
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -0.64%  '

# Row 13
$ws.Range("D13").Value = '37.91'
$ws.Range("E13").Value = '  -0.66%  '

# Row 14
$ws.Range("E14").Value = '  +0.88%  '

# Row 15
$ws.Range("D15").Value = '4.359.77'
$ws.Range("E15").Value = '  +0.93%  '

# Row 16
$ws.Range("D16").Value = '3.734.34'
$ws.Range("E16").Value = '  +0.99%  '

# Row 17
$ws.Range("D17").Value = '69.059.29'
$ws.Range("E17").Value = '  +2.10%  '

# Row 18
$ws.Range("D18").Value = '7.28'
$ws.Range("E18").Value = '  -0.15%  '

# Row 19
$ws.Range("E19").Value = '  -1.25%  '

# Row 20
$ws.Range("D20").Value = '17.13'
$ws.Range("E20").Value = '  -1.12%  '

# Row 21
$ws.Range("D21").Value = '10.70'
$ws.Range("E21").Value = '  +15.83%  '

# Row 22
$ws.Range("D22").Value = '491.70'
$ws.Range("E22").Value = '  +0.61%  '

# Row 23
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  -0.70%  '

# Row 24
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  +5.33%  '

# Row 25
$ws.Range("D25").Value = '84.70'
$ws.Range("E25").Value = '  -0.19%  '

# Row 26
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").Value = '12.27'
$ws.Range("E27").Value = '  -0.28%  '

# Row 28
$ws.Range("D28").Value = '10.10'
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("E29").Value = '  -0.31%  '

# Row 30
$ws.Range("E30").Value = '  +1.68%  '

# Row 31
$ws.Range("D31").Value = '2.49'
$ws.Range("E31").Value = '  +5.24%  '

# Row 32
$ws.Range("E32").Value = '  +3.94%  '

# Row 33
$ws.Range("D33").Value = '31.44'
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("D34").Value = '3.881.82'
$ws.Range("E34").Value = '  +1.15%  '

# Row 35
$ws.Range("E35").Value = '  -0.34%  '

# Row 36
$ws.Range("D36").Value = '3.671.55'
$ws.Range("E36").Value = '  +0.51%  '

# Row 37
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.41%  '

# Row 38
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  +1.64%  '

# Row 39
$ws.Range("D39").Value = '5.83'
$ws.Range("E39").Value = '  -0.25%  '

# Row 40
$ws.Range("E40").Value = '  +1.16%  '

# Row 41
$ws.Range("D41").Value = '0.323'
$ws.Range("E41").Value = '  -0.35%  '

# Row 42
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '430.97'
$ws.Range("E42").Value = '  +0.44%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.94'
$ws.Range("E43").Value = '  +3.35%  '

# Row 44
$ws.Range("D44").Value = '48.43'
$ws.Range("E44").Value = '  -1.30%  '

# Row 45
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  +0.76%  '

# Row 46
$ws.Range("E46").Value = '  +0.33%  '

# Row 47
$ws.Range("E47").Value = '  +0.03%  '

# Row 48
$ws.Range("D48").Value = '39.99'
$ws.Range("E48").Value = '  -1.53%  '

# Row 49
$ws.Range("D49").Value = '141.70'
$ws.Range("E49").Value = '  +0.47%  '

# Row 50
$ws.Range("D50").Value = '2.776.13'
$ws.Range("E50").Value = '  +0.89%  '

# Row 51
$ws.Range("E51").Value = '  +0.14%  '
